# Update the "Username" column (column E) values on Sheet1 to the new
# usernames, per the shared-strings diff:
#   JJANE        -> janeTheGem
#   john09       -> johntheGoat
#   adENA        -> lauraTheRainbow
#   sarat3ITA    -> saraThegold
#   linda11EAD1  -> lindaArrio
#   johnggDDS    -> joeBear

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "janeTheGem"
$ws.Range("E3").Value = "johntheGoat"
$ws.Range("E4").Value = "lauraTheRainbow"
$ws.Range("E5").Value = "saraThegold"
$ws.Range("E6").Value = "lindaArrio"
$ws.Range("E7").Value = "joeBear"
